$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1541.0588
$ws.Range("I132").Value = 1313.7142
$ws.Range("K132").Value = 3941.1426
$ws.Range("M132").Value = -1411.1426
$ws.Range("H137").Value = 31509.395
$ws.Range("I137").Value = 1164.5358
$ws.Range("J137").Value = 201440.6
$ws.Range("K137").Value = 3493.6074
$ws.Range("L137").Value = 604321.8
$ws.Range("M137").Value = -943.6074000000003
$ws.Range("N137").Value = -609421.8

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3584.9219
$ws.Range("I32").Value = 2730.7637
$ws.Range("K32").Value = 2730.7637
$ws.Range("M32").Value = -2443.7637
$ws.Range("H45").Value = 1608.1818
$ws.Range("I45").Value = 1099.5
$ws.Range("J45").Value = 1721.2222
$ws.Range("K45").Value = 1099.5
$ws.Range("L45").Value = 1721.2222
$ws.Range("M45").Value = -722.5
$ws.Range("N45").Value = -2475.2222
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("H132").Value = 2844.25
$ws.Range("I132").Value = 2984.2
$ws.Range("K132").Value = 8952.599999999999
$ws.Range("M132").Value = -6422.599999999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 120175.3
$ws.Range("I86").Value = 2448.5833
$ws.Range("K86").Value = 2448.5833
$ws.Range("M86").Value = -1325.5833
$ws.Range("H89").Value = 120175.3
$ws.Range("I89").Value = 2448.5833
$ws.Range("K89").Value = 12242.9165
$ws.Range("M89").Value = -6626.916499999999
$ws.Range("H105").Value = 2238.9048
$ws.Range("I105").Value = 2211.4211
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 2211.4211
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -464.4211
$ws.Range("N105").Value = -5994

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1836.5294
$ws.Range("I31").Value = 1263
$ws.Range("J31").Value = 2655.8572
$ws.Range("K31").Value = 1263
$ws.Range("L31").Value = 2655.8572
$ws.Range("M31").Value = -968
$ws.Range("N31").Value = -3245.8572
$ws.Range("H34").Value = 1836.5294
$ws.Range("I34").Value = 1263
$ws.Range("J34").Value = 2655.8572
$ws.Range("K34").Value = 1263
$ws.Range("L34").Value = 2655.8572
$ws.Range("M34").Value = -1061
$ws.Range("N34").Value = -3059.8572
$ws.Range("H58").Value = 1978327.2
$ws.Range("I58").Value = 3106969.8
$ws.Range("K58").Value = 3106969.8
$ws.Range("M58").Value = -3106766.8
$ws.Range("H86").Value = 2108.9
$ws.Range("I86").Value = 1596.6
$ws.Range("J86").Value = 2621.2
$ws.Range("K86").Value = 1596.6
$ws.Range("L86").Value = 2621.2
$ws.Range("M86").Value = -473.5999999999999
$ws.Range("N86").Value = -4867.2
$ws.Range("H89").Value = 2108.9
$ws.Range("I89").Value = 1596.6
$ws.Range("J89").Value = 2621.2
$ws.Range("K89").Value = 7983
$ws.Range("L89").Value = 13106
$ws.Range("M89").Value = -2367
$ws.Range("N89").Value = -24338
$ws.Range("H94").Value = 1137.3334
$ws.Range("I94").Value = 933.3333
$ws.Range("J94").Value = 1341.3334
$ws.Range("K94").Value = 933.3333
$ws.Range("L94").Value = 1341.3334
$ws.Range("M94").Value = -482.3333
$ws.Range("N94").Value = -2243.3334
$ws.Range("H99").Value = 1668935.5
$ws.Range("I99").Value = 5000599.5
$ws.Range("J99").Value = 3103.5
$ws.Range("K99").Value = 5000599.5
$ws.Range("L99").Value = 3103.5
$ws.Range("M99").Value = -4999101.5
$ws.Range("N99").Value = -6099.5
$ws.Range("H122").Value = 2727
$ws.Range("I122").Value = 1571.1428
$ws.Range("J122").Value = 27000
$ws.Range("K122").Value = 4713.428400000001
$ws.Range("L122").Value = 81000
$ws.Range("M122").Value = -2263.428400000001
$ws.Range("N122").Value = -85900
$ws.Range("H126").Value = 1668935.5
$ws.Range("I126").Value = 5000599.5
$ws.Range("J126").Value = 3103.5
$ws.Range("K126").Value = 15001798.5
$ws.Range("L126").Value = 9310.5
$ws.Range("M126").Value = -14999328.5
$ws.Range("N126").Value = -14250.5
$ws.Range("H132").Value = 2985.6667
$ws.Range("J132").Value = 3934.6667
$ws.Range("L132").Value = 11804.0001
$ws.Range("N132").Value = -16864.0001
$ws.Range("H134").Value = 1494.25
$ws.Range("J134").Value = 1289
$ws.Range("L134").Value = 3867
$ws.Range("N134").Value = -8937
$ws.Range("H135").Value = 60780
$ws.Range("J135").Value = 60780
$ws.Range("L135").Value = 60780
$ws.Range("N135").Value = -70920
$ws.Range("H136").Value = 1978327.2
$ws.Range("I136").Value = 3106969.8
$ws.Range("K136").Value = 9320909.399999999
$ws.Range("M136").Value = -9318359.399999999

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 543.3333
$ws.Range("I68").Value = 592
$ws.Range("K68").Value = 1776
$ws.Range("M68").Value = -965
$ws.Range("H71").Value = 543.3333
$ws.Range("I71").Value = 592
$ws.Range("K71").Value = 5328
$ws.Range("M71").Value = -1272
$ws.Range("H122").Value = 987.8333
$ws.Range("J122").Value = 1088.3334
$ws.Range("L122").Value = 9795.000599999999
$ws.Range("N122").Value = -14695.0006
$ws.Range("H132").Value = 950
$ws.Range("I132").Value = 950
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8550
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6020
$ws.Range("N132").Value = ""
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").Value = ""
$ws.Range("H134").Value = 2081.8333
$ws.Range("I134").Value = 1698.8
$ws.Range("K134").Value = 5096.4
$ws.Range("M134").Value = -26.39999999999964
$ws.Range("H137").Value = 2685.3684
$ws.Range("I137").Value = 1690
$ws.Range("J137").Value = 2950.8
$ws.Range("K137").Value = 5070
$ws.Range("L137").Value = 8852.400000000001
$ws.Range("M137").Value = 30
$ws.Range("N137").Value = -19052.4
$ws.Range("H138").Value = 3096.077
$ws.Range("I138").Value = 2449.9
$ws.Range("J138").Value = 5250
$ws.Range("K138").Value = 7349.700000000001
$ws.Range("L138").Value = 15750
$ws.Range("M138").Value = -2209.700000000001
$ws.Range("N138").Value = -26030
$ws.Range("H140").Value = 1423.5625
$ws.Range("I140").Value = 804.8570999999999
$ws.Range("J140").Value = 2604.7273
$ws.Range("K140").Value = 2414.5713
$ws.Range("L140").Value = 7814.1819
$ws.Range("M140").Value = 2765.4287
$ws.Range("N140").Value = -18174.1819
$ws.Range("H141").Value = 2714.7144
$ws.Range("I141").Value = 2714.7144
$ws.Range("K141").Value = 8144.1432
$ws.Range("M141").Value = -2964.1432

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""
$ws.Range("H113").Value = 1232
$ws.Range("I113").Value = 1098.5
$ws.Range("K113").Value = 1098.5
$ws.Range("M113").Value = 1071.5
$ws.Range("H122").Value = 2340
$ws.Range("I122").Value = 1733.3334
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 5200.0002
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -2750.0002
$ws.Range("N122").Value = -14650
$ws.Range("H126").Value = 2317585.8
$ws.Range("I126").Value = 2926803
$ws.Range("J126").Value = 2560
$ws.Range("K126").Value = 8780409
$ws.Range("L126").Value = 7680
$ws.Range("M126").Value = -8777939
$ws.Range("N126").Value = -12620

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3431.2222
$ws.Range("I7").Value = 2610.125
$ws.Range("K7").Value = 2610.125
$ws.Range("M7").Value = -2498.125
$ws.Range("H40").Value = 5063.875
$ws.Range("I40").Value = 1456.6364
$ws.Range("K40").Value = 1456.6364
$ws.Range("M40").Value = -1320.6364
$ws.Range("H61").Value = 1732.8889
$ws.Range("I61").Value = 1433
$ws.Range("K61").Value = 1433
$ws.Range("M61").Value = -1231
$ws.Range("H68").Value = 2334.7273
$ws.Range("I68").Value = 2438.9
$ws.Range("J68").Value = 1293
$ws.Range("K68").Value = 2438.9
$ws.Range("L68").Value = 1293
$ws.Range("M68").Value = -1689.9
$ws.Range("N68").Value = -2791
$ws.Range("H71").Value = 2334.7273
$ws.Range("I71").Value = 2438.9
$ws.Range("J71").Value = 1293
$ws.Range("K71").Value = 12194.5
$ws.Range("L71").Value = 6465
$ws.Range("M71").Value = -8450.5
$ws.Range("N71").Value = -13953
$ws.Range("H113").Value = 1732.8889
$ws.Range("I113").Value = 1433
$ws.Range("K113").Value = 1433
$ws.Range("M113").Value = 737
$ws.Range("H126").Value = 3431.2222
$ws.Range("I126").Value = 2610.125
$ws.Range("K126").Value = 7830.375
$ws.Range("M126").Value = -5360.375

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""
$ws.Range("H107").Value = 519.4483
$ws.Range("I107").Value = 425.42307
$ws.Range("K107").Value = 1276.26921
$ws.Range("M107").Value = 643.7307900000001
$ws.Range("H122").Value = 87825.336
$ws.Range("I122").Value = 112204
$ws.Range("K122").Value = 336612
$ws.Range("M122").Value = -334162
